$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.474.01'
$ws.Range("E2").Value = '  +2.18%  '
$ws.Range("D3").Value = '3.391.79'
$ws.Range("E3").Value = '  +1.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.63'
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.24'
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.88%  '
$ws.Range("E9").Value = '  +5.53%  '
$ws.Range("E10").Value = '  +1.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.67'
$ws.Range("E11").Value = '  +2.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000283'
$ws.Range("E12").Value = '  +3.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '681.49'
$ws.Range("E13").Value = '  -4.04%  '
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("D15").Value = '3.937.15'
$ws.Range("E15").Value = '  +1.51%  '
$ws.Range("D16").Value = '69.516.98'
$ws.Range("E16").Value = '  +2.23%  '
$ws.Range("D17").Value = '3.393.03'
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("E18").Value = '  +1.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.69'
$ws.Range("E19").Value = '  +0.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.31'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.906'
$ws.Range("E21").Value = '  +0.70%  '
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.22'
$ws.Range("E23").Value = '  +0.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '103.91'
$ws.Range("E24").Value = '  +3.47%  '
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("E26").Value = '  +0.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.78'
$ws.Range("E27").Value = '  +1.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '34.22'
$ws.Range("E28").Value = '  +2.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.74'
$ws.Range("E29").Value = '  +0.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.06'
$ws.Range("E30").Value = '  -0.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.19'
$ws.Range("E31").Value = '  +1.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '559.38'
$ws.Range("E32").Value = '  -2.49%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.107'
$ws.Range("E33").Value = '  +0.75%  '
$ws.Range("B34").Value = 'dogwifhat'
$ws.Range("C34").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.58'
$ws.Range("E34").Value = '  +3.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.28'
$ws.Range("E35").Value = '  +1.69%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").Value = '3.695.57'
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.61'
$ws.Range("E38").Value = '  +2.40%  '
$ws.Range("E39").Value = '  +4.82%  '
$ws.Range("E40").Value = '  +3.11%  '
$ws.Range("E41").Value = '  +1.69%  '
$ws.Range("D42").Value = '0.0₃0700'
$ws.Range("E42").Value = '  +2.67%  '
$ws.Range("E43").Value = '  +0.49%  '
$ws.Range("E44").Value = '  +3.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.31'
$ws.Range("E45").Value = '  -1.07%  '
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("E48").Value = '  +5.14%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.99'
$ws.Range("E50").Value = '  +1.51%  '
$ws.Range("E51").Value = '  +1.77%  '
